# Update the "department" column (C) on the "courses" sheet to split the
# single "FACULTY OF BUSINESS & TECHNOLOGY" label into more specific
# department / package names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("courses")

# row -> new department text
$departments = @{
    2  = "Business"
    3  = "Business"
    4  = "Business"
    5  = "Business"
    6  = "Business"
    7  = "Business"
    8  = "Business"
    9  = "Business"
    10 = "Business"
    11 = "Business"
    12 = "Business"
    13 = "Information Technology"
    14 = "Information Technology"
    15 = "Information Technology"
    16 = "Building and Construction"
    17 = "Packages"
    18 = "Packages"
    19 = "Packages"
    20 = "Packages"
    21 = "Packages"
    22 = "Packages"
}

foreach ($row in $departments.Keys) {
    $ws.Range("C$row").Value = $departments[$row]
}
